$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Build the two new cell-style "templates" that this workbook needs for
#    the new divider row (row 27): a border with a thin line on both the
#    top and the bottom (existing borders only had a thin bottom, or a
#    medium top+bottom). We derive them from the existing "thin bottom"
#    style (used on A4/B4 and A21/B21) so we keep reusing the same
#    underlying border/cellXf records instead of Excel inventing new,
#    unrelated ones.
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats
$tmpFont0 = $ws.Range("Z1")
$tmpFont0.Borders.Item(8).LineStyle = 1     # xlEdgeTop = thin continuous
$tmpFont0.WrapText = $true
$tmpFont0.HorizontalAlignment = -4142       # xlAlignNone (clear inherited "left")
$tmpFont0.VerticalAlignment = -4142         # xlAlignNone (clear inherited "top")

$tmpFont0.Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$tmpFont1 = $ws.Range("Z2")
$tmpFont1.Font.Size = 8

# Apply the two templates onto the real target cells of row 27.
$tmpFont0.Copy()
$ws.Range("A27:B27").PasteSpecial(-4122)
$tmpFont1.Copy()
$ws.Range("C27:E27").PasteSpecial(-4122)

$ws.Range("Z1:Z2").Clear()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Row 26 - the group that used to end the sheet now needs the
#    "thin bottom border" divider styling (same look as rows 4 and 21),
#    since a new group of rows (27-31) is appended after it.
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D4").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. New data rows 27-31.
# ---------------------------------------------------------------------------
$ws.Range("A27").Value = "SCRIPT/P01P04A/um1204.ssb"
$ws.Range("B27").Value = 274
$ws.Range("C27").Value = " I heard that the different kinds\nof [CS:I]Prize Tickets[CR] have something\nto do with your Explorer Rank…"
$ws.Range("D27").Value = " Говорят, что разные типы\n[CS:I]Призовых Билетов[CR] как-то связаны с Рангом\nИсследователя..."
$ws.Range("E27").Value = " Ãïâïñÿó, œóï ñàèîúå óéðú\n[CS:I]Ðñéèïâúö Áéìåóïâ[CR] ëàë-óï òâÿèàîú ò Ñàîãïí\nÉòòìåäïâàóåìÿ…"
$ws.Rows.Item(27).RowHeight = 43.2

$ws.Range("A28").Value = "SCRIPT/T01P02A/um1316.ssb"
$ws.Range("B28").Value = 255
$ws.Range("C28").Value = " [CS:N]Grovyle[CR], huh...?[K] Can't say I've\nheard of him, to be honest."
$ws.Range("D28").Value = " [CS:N]Гровайл[CR], да?..[K] Не могу сказать,\nчто я о нём что либо слышал."
$ws.Range("E28").Value = " [CS:N]Ãñïâàêì[CR], äà?..[K] Îå íïãô òëàèàóû,\nœóï ÿ ï îæí œóï ìéáï òìúšàì."
$ws.Rows.Item(28).RowHeight = 43.2

$ws.Range("A29").Value = "SCRIPT/T01P02A/um1402.ssb "
$ws.Rows.Item(29).RowHeight = 43.2

$ws.Range("A30").Value = "SCRIPT/T01P02A/um1608.ssb"
$ws.Range("B30").Value = 233
$ws.Range("C30").Value = " Doing nothing but waiting is\nfrustrating, but...[K]it's the great [CS:N]Dusknoir[CR]'s call."
$ws.Range("D30").Value = " Печально, что кроме ожидания\nмы ничего не можем сделать, но...[K] Таков\nзамысел великого [CS:N]Даскнуара[CR]."
$ws.Range("E30").Value = " Ðåœàìûîï, œóï ëñïíå ïçéäàîéÿ\níú îéœåãï îå íïçåí òäåìàóû, îï...[K] Óàëïâ\nèàíúòåì âåìéëïãï [CS:N]Äàòëîôàñà[CR]."
$ws.Rows.Item(30).RowHeight = 43.2

$ws.Range("B31").Value = 236
$ws.Range("C31").Value = " We have to trust him to do it."
$ws.Range("D31").Value = " Мы доверимся его словам."
$ws.Range("E31").Value = " Íú äïâåñéíòÿ åãï òìïâàí."

# ---------------------------------------------------------------------------
# 4. Apply the correct cell styles (font / border) to each new row, mirroring
#    the pattern already used for the earlier groups in the sheet:
#      - plain rows  -> style of A2/B2 (col A/B) + C2/D2/E2 (col C/D/E)
#      - divider row -> style of A21/B21 (col A/B) + C21/D21/E21 (col C/D/E)
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$ws.Range("B30").PasteSpecial(-4122)
$ws.Range("B31").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E31").PasteSpecial(-4122)

$ws.Range("A21").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("B21").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D21").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E21").Copy()
$ws.Range("E29").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5. Update the view so the new last cell is the active selection,
#    similar to how the author's own selection shifted down to E31.
# ---------------------------------------------------------------------------
$ws.Range("E31").Select()
